$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.34
$ws.Range("K2").Value = 7
$ws.Range("S2").Value = 2.48
$ws.Range("V2").Value = 1.07
$ws.Range("W2").Value = 3.95
$ws.Range("X2").Value = 26
$ws.Range("AB2").Value = 9.6
$ws.Range("AC2").Value = 15
$ws.Range("AE2").Value = 270
$ws.Range("AG2").Value = 11.5
